$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.867990493774414
$ws.Range("B1").Value = 3.485346794128418
$ws.Range("C1").Value = 1.841484546661377
$ws.Range("D1").Value = 1.455273628234863
$ws.Range("E1").Value = 1.331230878829956
